# New weekly record: insert a row before row 151, shifting the existing
# rows 151-177 down to 152-178, then populate the newly opened row 151
# with this week's data (same market/product/origin as the surrounding
# rows, new date, and the latest price readings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 151:177 down by one to open up a blank row at 151.
$ws.Rows.Item(151).Insert()

$ws.Cells.Item(151, 1).Value = 4
$ws.Cells.Item(151, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(151, 3).Value = "Los Lagos"
$ws.Cells.Item(151, 4).Value = 44522
$ws.Cells.Item(151, 5).Value = 10
$ws.Cells.Item(151, 6).Value = "Fruta"
$ws.Cells.Item(151, 7).Value = 100102
$ws.Cells.Item(151, 8).Value = "Cítricos"
$ws.Cells.Item(151, 9).Value = 100102006
$ws.Cells.Item(151, 10).Value = "Pomelo"
$ws.Cells.Item(151, 11).Value = "Start Ruby"
$ws.Cells.Item(151, 12).Value = "Primera"
$ws.Cells.Item(151, 13).Value = 60
$ws.Cells.Item(151, 14).Value = 11000
$ws.Cells.Item(151, 15).Value = 12000
$ws.Cells.Item(151, 16).Value = 11500
$ws.Cells.Item(151, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(151, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(151, 19).Value = 821
$ws.Cells.Item(151, 20).Value = 14
